$p = $ppt.ActivePresentation

# Slide 4: merge the title runs "Studio 2 Solution " + "(last class)" into one run.
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Studio 2 Solution (last class)"

# Slide 5: update the subtitle "Unit 2 - Class 4" -> "Unit 2 - Class 5"
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(2).TextFrame.TextRange.Text = "Unit 2 - Class 5"
